# "Tuesday4Feb25 time formatted good"
#
# The "Previsto" (scheduled time) column D had been filled with a single
# duplicated value for every stop (with one stray bad entry at D10) instead
# of the actual per-stop schedule. This fixes rows 3-14 with the correct,
# increasing h:mm times, and restores the sheet view (scrolled back to the
# top, selection parked on D15 - the next untouched row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

$ws.Range("D3").Value  = 0.3611111111111111
$ws.Range("D4").Value  = 0.37638888888888888
$ws.Range("D5").Value  = 0.37847222222222227
$ws.Range("D6").Value  = 0.38263888888888892
$ws.Range("D7").Value  = 0.38611111111111113
$ws.Range("D8").Value  = 0.3888888888888889
$ws.Range("D9").Value  = 0.39027777777777778
$ws.Range("D10").Value = 0.3923611111111111
$ws.Range("D11").Value = 0.39444444444444443
$ws.Range("D12").Value = 0.3972222222222222
$ws.Range("D13").Value = 0.39999999999999997
$ws.Range("D14").Value = 0.40277777777777773

# Move the selection (this also drops the old topLeftCell scroll offset,
# since the view is no longer pinned away from the top of the sheet).
$ws.Range("D15").Select()
